$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45/46 swap: EnergySwap moves up to row 45, PaxDollar moves down to row 46,
# each with refreshed Price/Volume figures.
$ws.Range('B45').Formula = 'EnergySwap'
$ws.Range('C45').Formula = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Formula = "'10.22"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Formula = '  -0.97%  '
$ws.Range('B46').Formula = 'PaxDollar'
$ws.Range('C46').Formula = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Formula = "'1.006"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Formula = '  +0.17%  '

# Price / Volume(1h) refreshes for all other rows
$ws.Range('D2').Formula = "'26.883.25"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Formula = '  -2.28%  '
$ws.Range('D3').Formula = "'1.834.39"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Formula = '  -1.72%  '
$ws.Range('D4').Formula = "'1.007"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Formula = '  +0.23%  '
$ws.Range('D5').Formula = "'310.43"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Formula = '  -1.80%  '
$ws.Range('E6').Formula = '  +0.19%  '
$ws.Range('D7').Formula = "'0.4615"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Formula = '  -1.46%  '
$ws.Range('D8').Formula = "'0.3667"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Formula = '  -1.69%  '
$ws.Range('D9').Formula = "'0.07176"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Formula = '  -2.87%  '
$ws.Range('D10').Formula = "'0.8805"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Formula = '  -1.05%  '
$ws.Range('D11').Formula = "'0.07864"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Formula = '  -0.96%  '
$ws.Range('D12').Formula = "'19.64"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Formula = '  -2.07%  '
$ws.Range('D13').Formula = "'1.846.94"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Formula = '  -1.65%  '
$ws.Range('D14').Formula = "'5.348"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Formula = '  -1.45%  '
$ws.Range('D15').Formula = "'6.405"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Formula = '  -3.03%  '
$ws.Range('D16').Formula = "'88.47"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Formula = '  -4.59%  '
$ws.Range('E17').Formula = '  +0.15%  '
$ws.Range('D18').Formula = "'0.000008758"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Formula = '  -2.01%  '
$ws.Range('D19').Formula = "'1.005"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Formula = '  +0.12%  '
$ws.Range('D20').Formula = "'26.923.59"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Formula = '  -2.23%  '
$ws.Range('E21').Formula = '  -2.78%  '
$ws.Range('D22').Formula = "'5.017"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Formula = '  -1.28%  '
$ws.Range('D24').Formula = "'1.983"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Formula = '  +5.75%  '
$ws.Range('D25').Formula = "'151.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Formula = '  -1.58%  '
$ws.Range('D26').Formula = "'18.25"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Formula = '  -1.39%  '
$ws.Range('D27').Formula = "'1.988"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Formula = '  -5.02%  '
$ws.Range('D28').Formula = "'113.85"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Formula = '  -2.53%  '
$ws.Range('E29').Formula = '  -4.27%  '
$ws.Range('D30').Formula = "'0.08848"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Formula = '  -0.72%  '
$ws.Range('D31').Formula = "'3.130"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Formula = '  +3.39%  '
$ws.Range('D32').Formula = "'0.7601"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Formula = '  +0.24%  '
$ws.Range('D33').Formula = "'4.470"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Formula = '  -0.41%  '
$ws.Range('E34').Formula = '  -3.21%  '
$ws.Range('D35').Formula = "'2.651"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Formula = '  +0.60%  '
$ws.Range('E36').Formula = '  +0.50%  '
$ws.Range('D37').Formula = "'0.01937"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Formula = '  -1.52%  '
$ws.Range('D38').Formula = "'2.937"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Formula = '  -1.74%  '
$ws.Range('D39').Formula = "'0.05144"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Formula = '  -2.49%  '
$ws.Range('D40').Formula = "'6.930"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Formula = '  -3.49%  '
$ws.Range('D41').Formula = "'0.4985"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Formula = '  -4.63%  '
$ws.Range('D42').Formula = "'0.1600"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Formula = '  -2.82%  '
$ws.Range('D43').Formula = "'8.314"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Formula = '  -0.90%  '
$ws.Range('D44').Formula = "'0.4704"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Formula = '  -3.62%  '
$ws.Range('D47').Formula = "'102.78"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Formula = '  -1.11%  '
$ws.Range('E48').Formula = '  -2.65%  '
$ws.Range('D49').Formula = "'0.06094"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Formula = '  -2.79%  '
$ws.Range('D50').Formula = "'64.92"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Formula = '  -1.65%  '
$ws.Range('D51').Formula = "'36.43"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Formula = '  -1.99%  '
